$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D29").Value = "[Self-Supervised Learning] Review: Bootstrap Your Own Latent -A New Approach to Self-Supervised Learning"
$ws.Range("E29").Value = "https://blog.promedius.ai/self-supervised-learning-review-bootstrap-your-own-latent-a-new-approach-to-self-supervised-learning/"

$ws.Range("D39").Value = "Deep Learning — Different Types of Autoencoders"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Deep-Learning-%E2%80%94-Different-Types-of-Autoencoders-1"

$ws.Range("D51").Value = "[MariaDB] enum 데이터 타입"
$ws.Range("E51").Value = "https://bskyvision.com/1152"

$ws.Range("D52").Value = "메타프로그래밍: 거울 테스트, 메타 인지, 스스로 발전하는 프로그램(aka 인공지능)"
$ws.Range("E52").Value = "http://ds.sumeun.org/?p=2323&utm_source=rss&utm_medium=rss&utm_campaign=%25eb%25a9%2594%25ed%2583%2580%25ed%2594%2584%25eb%25a1%259c%25ea%25b7%25b8%25eb%259e%2598%25eb%25b0%258d-%25ea%25b1%25b0%25ec%259a%25b8-%25ed%2585%258c%25ec%258a%25a4%25ed%258a%25b8-%25eb%25a9%2594%25ed%2583%2580-%25ec%259d%25b8%25ec%25a7%2580-%25ec%258a%25a4%25ec%258a%25a4%25eb%25a1%259c-%25eb%25b0%259c%25ec%25a0%2584-2"
